$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column from 2023-10-05 (45204) to 2023-10-08 (45207)
# for rows 2 through 5.
$ws.Range("C2:C5").Value = 45207
